# ----------------------------------------------------------------------
# plantilla-bitacora.docx - "Add some minor corrections to the template"
#
#   1) Insert a new empty, centered, bold paragraph mark right before the
#      "EVIDENCIA FOTOGRÁFICA:" heading paragraph (matching the pattern
#      of the other blank spacer paragraphs above it).
#   2) Drop the stale <w:lastRenderedPageBreak/> cached inside the
#      heading's run.
#   3) Insert a new empty, centered paragraph (lang=es-SV) right before
#      the "${evidencia_fotografia:720:480}" placeholder paragraph, and
#      fix the placeholder text itself to "${evidencia_fotografia}"
#      (dropping the bogus ":720:480" suffix and turning the old
#      gramStart/gramEnd proofing markers into spellStart/spellEnd
#      around the whole word).
#   4) Merge the three runs making up "${/evidencia}" into a single run.
# ----------------------------------------------------------------------

$d = $word.ActiveDocument
$wordmlNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ------------------------------------------------------------------
# (1) + (2) "EVIDENCIA FOTOGRÁFICA:" paragraph
# ------------------------------------------------------------------

# (2) First, rewrite the heading text in place so the stale
#     lastRenderedPageBreak cache gets dropped from the run.
$fixRng = $d.Content
$fixRng.Find.ClearFormatting()
$fixRng.Find.Replacement.ClearFormatting()
$fixRng.Find.Execute("FOTOGR?FICA:", $true, $false, $true, $false, $false, $true, 1, $false, `
    "FOTOGR" + [char]0xC1 + "FICA:", 2) | Out-Null

# (1) Locate the heading paragraph again and insert a new blank
#     centered/bold paragraph immediately before it.
$headRng = $d.Content
$headRng.Find.ClearFormatting()
$headRng.Find.Execute("EVIDENCIA FOTOGR*FICA:", $true, $false, $true, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headPara = $headRng.Paragraphs(1)
$insBefore = $d.Range($headPara.Range.Start, $headPara.Range.Start)
$insBefore.InsertXML("<w:p $wordmlNs><w:pPr><w:jc w:val='center'/><w:rPr><w:b/></w:rPr></w:pPr></w:p>")

# ------------------------------------------------------------------
# (3) "${evidencia_fotografia:720:480}" paragraph
# ------------------------------------------------------------------

# Insert a new blank centered paragraph (lang=es-SV) right before it.
$photoRng = $d.Content
$photoRng.Find.ClearFormatting()
$photoRng.Find.Execute('${evidencia_fotografia:720:480}', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$photoPara = $photoRng.Paragraphs(1)
$photoInsBefore = $d.Range($photoPara.Range.Start, $photoPara.Range.Start)
$photoInsBefore.InsertXML("<w:p $wordmlNs><w:pPr><w:jc w:val='center'/><w:rPr><w:lang w:val='es-SV'/></w:rPr></w:pPr></w:p>")

# Re-find the placeholder paragraph (its position shifted after the
# insert above) and rewrite its runs / proofErr markers.
$photoRng2 = $d.Content
$photoRng2.Find.ClearFormatting()
$photoRng2.Find.Execute('${evidencia_fotografia:720:480}', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$photoPara2 = $photoRng2.Paragraphs(1)
$photoXml = "<w:p $wordmlNs>" + `
    "<w:pPr><w:jc w:val='center'/><w:rPr><w:lang w:val='es-SV'/></w:rPr></w:pPr>" + `
    "<w:r><w:rPr><w:lang w:val='es-SV'/></w:rPr><w:t>`${</w:t></w:r>" + `
    "<w:proofErr w:type='spellStart'/>" + `
    "<w:r><w:rPr><w:lang w:val='es-SV'/></w:rPr><w:t>evidencia_fotografia</w:t></w:r>" + `
    "<w:proofErr w:type='spellEnd'/>" + `
    "<w:r><w:rPr><w:lang w:val='es-SV'/></w:rPr><w:t>}</w:t></w:r>" + `
    "</w:p>"
$photoPara2.Range.InsertXML($photoXml)

# ------------------------------------------------------------------
# (4) "${/evidencia}" paragraph - merge the 3 runs into 1
# ------------------------------------------------------------------

$closeRng = $d.Content
$closeRng.Find.ClearFormatting()
$closeRng.Find.Execute('${/evidencia}', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$closePara = $closeRng.Paragraphs(1)
$closeXml = "<w:p $wordmlNs>" + `
    "<w:pPr><w:jc w:val='center'/><w:rPr><w:lang w:val='es-SV'/></w:rPr></w:pPr>" + `
    "<w:r><w:rPr><w:lang w:val='es-SV'/></w:rPr><w:t>`${/evidencia}</w:t></w:r>" + `
    "</w:p>"
$closePara.Range.InsertXML($closeXml)

$d.Save()
